# Edit the "Q1" sheet (first sheet, sheet1.xml):
#  - swap the "Code" / "Convention" columns (B <-> C) for header + all data rows
#  - add Grouping ("B"/"L") + GroupOrder values that were missing for a few rows
#  - move the bestFit column width from column B to column C
#  - change the active selection to F1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Q1")

# --- Row 1: header swap (Convention/Code -> Code/Convention) ---
$ws.Cells.Item(1, 2).Value = "Code"
$ws.Cells.Item(1, 3).Value = "Convention"

# --- Data rows 2-14: column B held the ICD-10 code, column C held the fixed
#     "ACS1001" convention value. After the edit, B holds the code and C
#     holds "ACS1001" (i.e. B and C are swapped). ---
$codes = @{
    2  = "I25.11"
    3  = "I25.12"
    4  = "I10"
    5  = "J98.5"
    6  = "T81.2"
    7  = "S25.0"
    8  = "Y60.0"
    9  = "Y92.22"
    10 = "J96.09"
    11 = "I48.9"
    12 = "J98.1"
    13 = "E87.6"
    14 = "F17.1"
}

foreach ($r in 2..14) {
    $ws.Cells.Item($r, 2).Value = $codes[$r]
    $ws.Cells.Item($r, 3).Value = "ACS1001"
}

# --- New Grouping / GroupOrder values ---
$ws.Cells.Item(4, 4).Value = "B"
$ws.Cells.Item(4, 6).Value = 1

$ws.Cells.Item(5, 4).Value = "B"

$ws.Cells.Item(6, 6).Value = 2

$ws.Cells.Item(13, 4).Value = "L"
$ws.Cells.Item(14, 4).Value = "L"

# --- Column widths: the bestFit width that used to belong to column B now
#     belongs to column C (column B goes back to an unset/default width). ---
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(2).ColumnWidth
$ws.Columns.Item(2).ColumnWidth = 8.43

# --- Selection moves to F1 ---
$ws.Range("F1").Select()
